$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Debatten om Danmarks indvandrings- og integrationspolitik fokuserer på kravene til sprogkundskaber og indfødsret, der skaber politisk splittelse og behov for strammere regler for at beskytte dansk identitet.'
$ws.Range("B3").Value = 'Diskussionen om diskrimination og lighed centrerer sig om indvandringspolitikken og uenigheder om, hvordan man skal forholde sig til indfødsret og permanent opholdstilladelse i forhold til danske borgere og fremmede.'
$ws.Range("B4").Value = 'Der er behov for strammere betingelser og mere gennemsigtighed i dansk indfødsret for at sikre samfundets sammenhængskraft og håndtere udfordringer relateret til indvandrere, samtidig med at ansøgere, der opfylder kravene, bør have ret til dansk statsborgerskab.'
$ws.Range("B5").Value = 'Høje boglige krav kan være en hindring for indvandrere, der ønsker selvforsørgelse, og lighed for loven samt forbud mod diskrimination er centrale borgerlige frihedsrettigheder.'
$ws.Range("B6").Value = 'Der er en anbefaling om at godkende lovforslaget om dansk statsborgerskab for langvarige beboere, med fokus på integration og sprogkrav, samtidig med kritik af Dansk Folkepartis syn på indfødsret og integration.'
$ws.Range("B7").Value = 'Der stilles spørgsmål ved det danske tankesæt og Dansk Folkepartis holdning til dispensation fra sprogkravet for indfødsret.'
$ws.Range("B8").Value = 'Indfødsret er afgørende for at bevare dansk identitet og kultur, og kravene til statsborgerskab bør skærpes for at sikre nye borgeres bidrag til fællesskabet.'
$ws.Range("B9").Value = 'Der er alvorlige problemer med retssikkerheden i tildelingen af dansk statsborgerskab, præget af afvisninger uden begrundelse og stramninger, som kan føre til eksklusion og radikalisering af ansøgere, mens der er behov for mere ensartet behandling og fokus på at hjælpe udsatte grupper.'
$ws.Range("B10").Value = 'Danmark oplever en folkevandring med mange ansøgere til statsborgerskab, hvor der forventes tilpasning til danske værdier og krav om sprogkundskaber.'
$ws.Range("B11").Value = 'Lovforslaget om indfødsret i Danmark belyser kompleksiteten og uretfærdigheden i lovgivningen, især for sårbare grupper som torturofre og asylbørn, samtidig med at der fremsættes ændringsforslag for at rette fejl i registreringen af statsborgerskab.'
$ws.Range("B12").Value = 'Der er en vilje til at finde løsninger for at give børn dansk statsborgerskab, især i tilfælde af fejlregistreringer eller manglende opmærksomhed.'
$ws.Range("B13").Value = 'Der er uenighed om indfødsret, hvor taleren mener, at sygdom ikke kan legitimere statsborgerskab, og at anerkendelse af folk og nation er vigtig.'
$ws.Range("B14").Value = 'Det Radikale Venstre støtter nye danske statsborgere, men ønsker ændringer i indfødsretsaftalen og kravene til statsborgerskab for at fremme integrationen.'
$ws.Range("B15").Value = 'Der opfordres til at byde nye statsborgere velkommen og anerkende deres engagement i det danske samfund, samtidig med at der kritiseres for uretfærdige begrænsninger i statsborgerskabsprocessen, særligt for dem med psykiske lidelser.'
$ws.Range("B16").Value = 'Den konservative folketingsgruppe støtter rimelige krav og særlige dispensationer for dansk statsborgerskab, men kritiserer den nye indfødsretsaftale for ikke at adressere vigtige problemstillinger.'
$ws.Range("B17").Value = 'Der er behov for klare krav til dansk sprogkundskab for statsborgerskab, og muligheden for dispensation ved posttraumatisk stress-syndrom skal vurderes grundigere.'
$ws.Range("B18").Value = 'Lovforslaget vil give cirka 1.800 personer, herunder over 1.000 børn, dansk statsborgerskab, men der er problemer med indfødsretsprøven, som Socialdemokraterne ønsker at løse.'
$ws.Range("B19").Value = 'Der er bekymring for en uretfærdig håndtering af indfødsretsprøven, og der efterlyses seriøst arbejde i stedet for hurtige løsninger.'
$ws.Range("B20").Value = 'Statsborgerskabsdagen fejrer nye danskere, men der er bekymring og frustration over strenge regler, der begrænser tildelingen af dansk statsborgerskab og integration.'
$ws.Range("B21").Value = 'Lovforslaget om indfødsret godkendes med fokus på dansk kultur og historie, men der er bekymringer om tidligere indhold og afvisning af dobbelt statsborgerskab.'
$ws.Range("B22").Value = 'Repræsentanten lykønsker kommende statsborgere og fremhæver vigtigheden af integration, mens der kritiseres de skærpede krav for statsborgerskab.'
$ws.Range("B23").Value = 'En betydningsfuld dag for 2.315 voksne og 810 børn, der opnår dansk statsborgerskab, hvilket understreger vigtigheden af integration og respekt for deres identitetsvalg.'
$ws.Range("B24").Value = 'Der er delte meninger om kravene til sprogprøven og vurderingsprocessen for indfødsret, hvor nogle mener, at kravene er rimelige, mens andre kalder på mere fleksibilitet for at fremme integration.'
$ws.Range("B25").Value = 'Både Morten Østergaard og Anne-Marie Meldgaard understreger vigtigheden af indfødsret og de konsekvenser, det har for ansøgerne.'
$ws.Range("B26").Value = 'Der er bekymringer og håb omkring indfødsret for personer med arabiske navne, samt en debat om kravene til statsborgerskab og integration i det danske samfund.'
$ws.Range("B27").Value = 'Ministerens erkendelse af vanskelighederne ved at opnå dansk statsborgerskab står i kontrast til regeringens beslutning om at skærpe integrationsreglerne.'
$ws.Range("B28").Value = 'Diskussionen om dansk statsborgerskab fokuserer på kravene til sprogkundskaber, viden om Danmark og ordentlig livsførelse, samt betydningen af indfødsret og konsekvenserne af dobbelt statsborgerskab.'
$ws.Range("B29").Value = 'Lovforslaget L 191 omhandler tildeling af dansk statsborgerskab til borgere, der ønsker at blive en del af det danske samfund, samtidig med at der er misforståelser omkring dobbelt statsborgerskab og dets konsekvenser for dansk identitet.'
$ws.Range("B30").Value = 'Den 26. april fejrede 1.500 nytilkomne danske statsborgere deres opnåede statsborgerskab, hvilket understreger vigtigheden af demokratisk deltagelse og rettigheder.'
$ws.Range("B31").Value = 'Dansk statsborgerskab giver tryghed til usikre personer, men der er misforståelser om kravene, herunder opgivelse af tidligere statsborgerskaber, og danske embedsmænd følger lovgivningen uden at tage hensyn til følelser i ansøgningsprocessen.'
$ws.Range("B32").Value = 'Der apporteres om behovet for Folketingets selvbestemmelse over indfødsret og afvisning af ministerens påstand om kommunistlovgivningens oprindelse.'
$ws.Range("B33").Value = 'Debatten om indfødsret i Danmark fokuserer på konflikten mellem opfyldelse af internationale konventioner og beskyttelse af grundloven samt Folketingets suverænitet i beslutninger om statsborgerskab.'
$ws.Range("B34").Value = 'Der er en fælles støtte til lovforslaget om dansk indfødsret, men samtidig en kritik af kravene og praksis omkring statsborgerskab, herunder ønsket om at tillade dobbelt statsborgerskab og bevare Folketingets autoritet i beslutningsprocessen.'
$ws.Range("B35").Value = 'Teksten omhandler tildeling af dansk statsborgerskab, vigtigheden af troskab over for danske værdier, og de krav og forberedelser, der er nødvendige for at blive statsborger i Danmark.'
$ws.Range("B36").Value = 'Der er bekymringer om demokratisk praksis og retssikkerhed i forbindelse med lovforslaget om statsborgerskab, herunder eksklusion af personer, PET''s integritet og ministerens manglende svar på vigtige spørgsmål.'
$ws.Range("B37").Value = 'Debatten i Folketinget omhandler balancen mellem Danmarks sikkerhed og overholdelsen af internationale konventioner, hvor der er delte meninger om håndtering af statsborgerskab for personer, der anses for sikkerhedstrusler.'
$ws.Range("B38").Value = 'Der er enighed om den nuværende retstilstand vedrørende efterretningstjenesten, men bekymringer om håndtering af sikkerhedsrisici og behovet for ændringer i konventioner diskuteres.'
$ws.Range("B39").Value = 'Der er enighed om nødvendigheden af at reformere statsborgerskabs krav, sikre integration af nye danskere og følge Politiets Efterretningstjenestes anbefalinger for at beskytte rigets sikkerhed, samtidig med at der er bekymringer om håndteringen af statsborgerskab for kriminelle.'
$ws.Range("B40").Value = 'SF og Socialdemokratiet byder nye statsborgere velkommen, men rejser spørgsmål om sagsbehandling og ønsket om klarhed over procedurerne.'
$ws.Range("B41").Value = 'Dansk Folkeparti afviser loven om indfødsrets meddelelse, da de mener, den krænker grundloven og underminerer Folketingets myndighed.'
$ws.Range("B42").Value = 'Konservative støtter lovforslaget om dansk statsborgerskab for dem, der opfylder kravene, og fremhæver vigtigheden af at anerkende ansøgeres bestræbelser og sikre politisk opbakning.'
$ws.Range("B43").Value = 'Lovforslaget om automatisk tildeling af statsborgerskab til statsløse personer er både kritiseret og bifaldet, hvilket rejser spørgsmål om retfærdighed, sikkerhed og behandling af personer født i Danmark.'
$ws.Range("B44").Value = 'Statsborgerskab tildeles 1.157 voksne og 487 børn i Danmark, hvilket anerkender deres tilknytning til landet og sikrer deres danske identitet.'
$ws.Range("B45").Value = 'Der er en debat om tildeling af dansk statsborgerskab til udlændinge, især statsløse, hvor bekymringer om sikkerhed og krav til ansøgere står centralt, samtidig med at der er anerkendelse af betydningen af aktiv deltagelse i civilsamfundet.'
$ws.Range("B46").Value = 'Der er en fælles stræben efter at forbedre adgangen til dansk statsborgerskab og udrydde statsløshed, samtidig med at forskellige politiske partier har forskellige holdninger til indfødsret og krav til medlemskab.'
$ws.Range("B47").Value = 'Regeringen og Folketinget kritiseres for at tildele statsborgerskab til kriminelle statsløse uden samme krav som andre ansøgere, hvilket skaber bekymringer om sikkerhed og moral.'
$ws.Range("B48").Value = 'Der er bekymring over tildelingen af dansk statsborgerskab til personer, vurderet som sikkerhedsrisikoer, som rejser spørgsmål om internationale konventioners fortolkning og Danmarks integration og statsborgerskabspolitik.'
$ws.Range("B49").Value = 'Dansk Folkeparti afviser lovforslaget om statsborgerskab med fokus på integration, ligebehandling og krav til ansøgere, samtidig med at de kritiserer lempelser i danskkundskaber og relationer til statsløse konventioner.'
$ws.Range("B50").Value = 'Dansk Folkepartis politik rejser spørgsmål om diskrimination mod muslimer og parallelsamfund i Danmark, samt uklarheder i deres krav til statsborgerskab, hvilket skaber bekymring for demokratiske værdier.'
$ws.Range("B51").Value = 'Der er bred politisk støtte til et lovforslag om dansk statsborgerskab for 1.109 personer, men samtidig er der bekymringer om kravene og ønsker om strengere regler.'
$ws.Range("B52").Value = 'Der er en fælles bekymring for indvandringens indvirkning på dansk kultur og samhørighed, samt en opfordring til at begrænse indvandringen fra ikkevestlige lande og stille krav til statsborgerskab for at sikre integration.'
$ws.Range("B53").Value = 'Debatten om dansk statsborgerskab fremhæver både glæden ved indfødsret for mange borgere og bekymringer over tildele statsborgerskab til personer vurderet som sikkerhedsrisikoer, hvilket rejser spørgsmål om ansvar, sikkerhed og politiske beslutninger.'
$ws.Range("B54").Value = 'Der er en debat om offentliggørelse af fortrolige oplysninger og håndtering af statsborgerskab, hvor forskellige synspunkter inden for Socialdemokratiet og spørgsmål om lovgivningens rammer spiller en central rolle.'
$ws.Range("B55").Value = 'Debatten om lovforslaget om statsborgerskab viser en splittelse mellem partierne, hvor fokus er på procedure og jura, samtidig med at der er bekymringer om behandlingen af fortrolige oplysninger og behandlingen af ansøgninger.'
$ws.Range("B56").Value = 'Der er en debat om statsløse personer, flygtninge- og integrationslovgivning samt tildeling af statsborgerskab i Danmark, hvor regeringens nuværende politik og tidligere fejl, herunder sagsbehandlingstider og PET''s rolle, drøftes.'
$ws.Range("B57").Value = 'Der er blevet givet mange dispensationer til dansk statsborgerskab, hvilket har ført til lempelser af kravene, men Dansk Folkeparti ønsker strammere krav.'
$ws.Range("B58").Value = 'Tildelingen af dansk statsborgerskab bør fokusere mere på integration og sprogkundskaber for at sikre, at nyankomne bedre kan blive en del af det danske samfund.'
$ws.Range("B59").Value = 'Lovforslaget skaber bekymring ved at blande berettigede statsborgerskabsansøgere med dem, der får dispensation, men Liberal Alliance vælger at støtte det på grund af overholdelse af retningslinjerne.'
$ws.Range("B60").Value = 'Lovforslaget om dansk statsborgerskab får positiv støtte fra konservative, der værdsætter statsborgerskabsdagen som en vigtig tradition.'
$ws.Range("B61").Value = '679 voksne og ca. 425 børn opnår dansk statsborgerskab, hvilket indebærer opfyldelse af betingelser og anerkendelse af deres unikke baggrunde og udfordringer.'
$ws.Range("B62").Value = 'Dagen fejres af 679 nye danske statsborgere, men processen for at opnå statsborgerskab er lang og kompliceret trods forbedringer for særlige grupper.'
$ws.Range("B63").Value = 'Velkomsten af nye statsborgere til det danske fællesskab og anerkendelsen af muligheden for statsborgerskab i flere lande er vigtige og glædelige begivenheder.'
$ws.Range("B64").Value = 'Debatten om indfødsretslovforslaget centrerer sig om principielle spørgsmål ved tildeling af statsborgerskab, hvor der er fokus på politisk anstændighed, integration af ansøgere, samt nødvendigheden af at respektere tidligere aftaler og sikre ensartet behandling.'
$ws.Range("B65").Value = 'Der er behov for en grundig vurdering af ansøgere om dansk statsborgerskab, der inkluderer loyalitet over for demokratiet og skelnen mellem forskellige holdninger samt en politisk beslutningsproces fremfor medicinske vurderinger.'
$ws.Range("B66").Value = 'Debatten om dansk statsborgerskab centrerer sig om krav til sprog, selvforsørgelse, politiske ytringer og retssikkerhed, hvor der er bekymringer om, hvordan disse faktorer påvirker tildelingen af statsborgerskab og integration i samfundet.'
$ws.Range("B67").Value = 'Fru Lotte Rod nævnes gentagne gange i teksten.'
$ws.Range("B68").Value = 'Dansk Folkeparti mener, at tildelingen af dansk statsborgerskab bør prioritere Danmarks fremtid og sammenhængskraft.'
$ws.Range("B69").Value = 'Der er enighed om, at unødvendige forhindringer i indfødsretssager skader samfundets sammenhængskraft, mens det er positivt at bevilge statsborgerskab til mange individer.'
$ws.Range("B70").Value = 'Der er behov for alvorlige ændringer i sagsbehandlingen af ansøgninger, herunder tydeligere retningslinjer, for at sikre kvalificerede beslutninger om integration og tage hensyn til menneskers liv.'
$ws.Range("B71").Value = 'Der er afstemninger om forskellige ændringsforslag, hvor flertallet støtter flere forslag, mens nogle enkelte forslag bliver forkastet.'
$ws.Range("B72").Value = 'Der er stolthed over nye danske statsborgere og anerkendelse af deres udfordringer, samtidig med bekymringer om retssikkerhed og krav til integration, herunder sprog og selvforsørgelse.'
$ws.Range("B73").Value = 'Tillykke til de nye danskere, men bekymringer om retssikkerhed for udelukkede ansøgere er rejst.'
$ws.Range("B74").Value = 'At blive dansk statsborger kræver opfyldelse af strenge krav og anses som et stort privilegium, hvilket fremhæves af tildelingen af statsborgerskab til mange voksne og børn.'
$ws.Range("B75").Value = 'Lovforslaget L 71 muliggør, at 3.185 nye danskere bliver statsborgere, hvilket anerkender deres bidrag og opfordrer til aktiv samfundsdeltagelse.'
$ws.Range("B76").Value = 'Der er enighed om at fastsætte kriterier for at reducere antallet af personer, der opnår dansk statsborgerskab, samtidig med at der diskuteres politisk stabilitet og sagsbehandlingstider i asyl- og statsborgerskabsprocessen.'
$ws.Range("B77").Value = 'Debatten om dansk statsborgerskab fokuserer på at skærpe betingelserne og anerkende nye danskeres ønske om integration, samtidig med at der er bekymringer om indvandringens indflydelse på Danmarks fremtid.'
$ws.Range("B78").Value = 'Lovforslaget om tildeling af dansk statsborgerskab til 2.139 borgere fremhæver vigtigheden af aktiv samfundsdeltagelse, respekt for danske værdier og anerkendelse af nye statsborgere som ligestillede medlemmer af det danske samfund.'
$ws.Range("B80").Value = 'Der er behov for et mere retfærdigt og gennemsigtigt system for sagsbehandlingen af ansøgninger om dansk statsborgerskab, da mange efterkommere stadig føler sig afvist i samfundet.'
$ws.Range("B81").Value = 'Lovforslaget om dansk statsborgerskab understreger vigtigheden af aktiv deltagelse i demokratiet og anerkender grundlæggende værdier, samtidig med at der er politisk støtte og bekymringer om sagsbehandlingstider.'
$ws.Range("B82").Value = 'Der er bekymring om integrationen af nye danske statsborgere, hvilket har ført til krav om strammere regler og en debat om, hvordan statsborgerskab tildeles.'
$ws.Range("B83").Value = 'Alternativet byder de nye danske statsborgere velkommen og opfordrer til aktiv deltagelse i samfundet.'
$ws.Range("B84").Value = 'Statsborgerskab er essentielt for integration og anerkendelse af individer i Danmark, da det muliggør fuld deltagelse i samfundet og styrker tilhørsforholdet.'
$ws.Range("B85").Value = 'Der er en ambivalens omkring tildeling af dansk statsborgerskab, hvor der ønskes strammere betingelser for at sikre integration og støtte til dem, der respekterer danske værdier, samtidig med at man undgår at inkludere dem, der modsætter sig demokratiet.'
$ws.Range("B86").Value = 'Venstre anerkender betydningen af statsborgerskab for voksne og børn, men der er bekymringer om partiets tilgang til dobbelt statsborgerskab og loyalitet over for flere lande.'
$ws.Range("B87").Value = 'Der er bekymring om indfødsret i Danmark, som bør tildeles under strenge betingelser for at sikre bedre integration og samfundsloyalitet, samtidig med at der anerkendes indsatsen fra nye statsborgere.'
$ws.Range("B88").Value = 'Der er behov for at lette adgangen til statsborgerskab i Danmark, da nuværende strenge krav og bureaukrati gør det svært for mange, især sårbare grupper, at opnå dette.'
$ws.Range("B89").Value = 'Debatten om statsborgerskab i Danmark fokuserer på krav om dansk identitet, sprogbeherskelse og samfundsengagement, samt nødvendigheden af klare kriterier og screeningsprocesser for ansøgere.'
$ws.Range("B90").Value = 'Der er en fælles bekymring for, hvordan udlændingeloven og statsborgerskabsforslag påvirker Danmarks sikkerhed og samfundsmæssige bidrag fra udlændinge.'
$ws.Range("B91").Value = 'Der er en bred debat om tildeling af dansk statsborgerskab, hvor flere partier anerkender glæden ved nye statsborgere, men kritiserer strenge regler og politiske beslutninger, mens andre ønsker at beskytte statsborgerskabet som en værdifuld gave.'
$ws.Range("B92").Value = 'Nye danskere anerkendes for deres bidrag, mens der stilles krav om forpligtelse til det danske samfund, samtidig med politisk debat om statsborgerskab og stramning af regler mod kriminalitet.'
$ws.Range("B93").Value = 'Debatten om dansk statsborgerskab centrerer sig om kravene til ansøgere, retssikkerhed, integration, og politiske værdier, med fokus på at sikre, at kun velintegrerede personer får statsborgerskab, mens der er bekymringer over ændringer i reglerne og politisk samarbejde.'
$ws.Range("B94").Value = 'Der er en politisk debat om indfødsret og statsborgerskab i Danmark, hvor partierne udtrykker bekymring for integration, krav til ansøgere og prioritering af danske værdier, især i relation til personer med muslimsk baggrund.'
$ws.Range("B95").Value = 'Der er en fælles anerkendelse af behovet for at revidere reglerne for dansk statsborgerskab, med fokus på retfærdighed, lighed og respekt for ansøgeres rettigheder, samtidig med at der tages hensyn til sundhedsanbefalinger.'
$ws.Range("B96").Value = 'Debatten om dansk statsborgerskab fokuserer på kravene til ansøgere, herunder håndtryk og kriminel baggrund, samt den politiske konflikt mellem partierne om udlændingepolitikken og tildeling af statsborgerskab.'
$ws.Range("B97").Value = 'Tillykke til nye danskere med statsborgerskab, hvor fokus er på samfundsengagement og tilpasning, samtidig med bekymringer om lovgivning, herunder håndtrykskravet og indvandreres kriminalitet.'
$ws.Range("B98").Value = 'Debatten om indfødsret i Danmark fokuserer på kravene til statsborgerskab, fortrolighed i udvalgsarbejde, og hvordan forskellige grupper opfattes i forhold til tildeling af dansk statsborgerskab.'
$ws.Range("B99").Value = 'Lovforslaget om massetildeling af danske statsborgerskaber kritiseres for at blande kvalificerede ansøgere med mindre kvalificerede, hvilket rejser behovet for strengere udvælgelse og klare regler for statsborgerskab.'
$ws.Range("B100").Value = 'Debatten om tildeling af dansk statsborgerskab er præget af utilfredshed med ministerens manglende gennemsigtighed og samarbejde, hvilket skaber bekymringer om retfærdighed, strammere krav og mulige konsekvenser for personer med kriminel baggrund.'
$ws.Range("B101").Value = 'Debatten om dansk statsborgerskab fokuserer på krav som sprogkundskaber, lovlydighed og selvforsørgelse samt bekymringer om kriminalitet og integration, hvilket fører til forslag om strammere regler og reformer i behandlingen af ansøgninger.'
$ws.Range("B102").Value = 'Statsborgerskab i Danmark kræver klare kriterier og en retfærdig proces, samtidig med at der er politiske uenigheder om behandling af ansøgere fra forskellige lande.'
$ws.Range("B103").Value = 'Der er behov for reform af reglerne for dansk statsborgerskab for at inkludere flere, mens debatten om indfødsretslovforslaget omhandler både tildeling af statsborgerskab og nye skærpede krav.'
$ws.Range("B104").Value = 'Der er blevet stemt om forskellige ændringsforslag til et lovforslag, hvor et flertal har vedtaget flere ændringer, mens der også er rejst bekymringer om lovgivningsprocessen og grundlovens bestemmelser vedrørende statsborgerskab og skat.'
$ws.Range("B105").Value = 'Der er behov for at anerkende og byde velkommen til borgere, der opfylder betingelserne for dansk statsborgerskab, samtidig med at der stilles krav om respekt for danske værdier og en differentieret tilgang til ansøgere baseret på oprindelsesland.'
$ws.Range("B106").Value = 'Der er en fælles bekymring blandt partierne om ansvarligheden ved tildeling af danske statsborgerskaber og betydningen af, at ansøgerne bidrager positivt til samfundet.'
$ws.Range("B107").Value = 'Der er en bred enighed blandt partierne om behovet for at forbedre retssikkerheden og proces for dansk statsborgerskab, samtidig med at der understreges vigtigheden af strenge kriterier og engagement i det danske samfund.'
$ws.Range("B108").Value = 'Debatten om dansk statsborgerskab fokuserer på krav til integration, håndtering af grov kriminalitet og beskyttelse af ofre, samtidig med at der stilles spørgsmål til de forskellige politiske partiers holdninger til statsborgerskabsfratagelse og internationale konventioner.'
$ws.Range("B109").Value = 'Debatten om dansk statsborgerskab fokuserer på behovet for klare kriterier og individuel vurdering af ansøgere, samtidig med at der udtrykkes bekymring for samfundets sammenhængskraft og retfærdighed i tildelingsprocessen.'
$ws.Range("B110").Value = 'Diskussionen om indfødsret fokuserer på retfærdighed og retssikkerhed for alle borgere, inklusive hvordan tidligere kriminelle kan rehabiliteres og genopnå statsborgerskab, samt vigtigheden af at sikre, at personer med ekstremistiske holdninger ikke får indfødsret.'
$ws.Range("B111").Value = 'Der er bekymringer om den nye indfødsretsaftale, da den kan føre til diskrimination og underkende vigtigheden af mangfoldighed og anerkendelse mellem mennesker.'
$ws.Range("B112").Value = 'Kritikken af massetildelingen af statsborgerskaber fokuserer på manglende kriterier og uigennemsigtighed, hvilket fører til tildeling af statsborgerskab til uegnede personer.'
$ws.Range("B113").Value = 'Der er en fælles anerkendelse af betydningen af statsborgerskab som en fundamental del af demokratiet, med fokus på ligebehandling, krav om integration og forskellige syn på reglerne for indfødsret.'
$ws.Range("B114").Value = 'Statsborgerskab ses som et privilegium, der bør tildeles værdige ansøgere, hvilket skaber debat om antallet af ansøgere, især fra MENAPT-lande, mens fejringen af de, der opnår dansk statsborgerskab, markerer berigelsen af det danske fællesskab.'
$ws.Range("B115").Value = 'Lovforslaget om tildeling af dansk statsborgerskab mødes med kritik for manglende kontrol og sagsbehandling, hvilket rejser bekymringer om sikkerhed, ansvarlighed og integration af både kriminelle og lovlydige individer.'
$ws.Range("B116").Value = 'Statsborgerskab i Danmark fremkalder både glæde og bekymringer over integration, værdier og sikkerhed, med opfordringer til en mere restriktiv tilgang til tildeling af statsborgerskaber, især for personer med problematiske holdninger.'
$ws.Range("B117").Value = 'Tildeling af dansk statsborgerskab kræver opfyldelse af strenge krav og rejser vigtige spørgsmål om integration, demografiske ændringer og samfundets sammenhængskraft, samtidig med at der er behov for retfærdige regler og anerkendelse af dem, der fortjener det.'
$ws.Range("B118").Value = 'Flere ændringsforslag er blevet afstemte og forkastet, mens lovforslaget til sidst er blevet vedtaget og sendt til statsministeren.'
$ws.Range("B119").Value = 'Alle ændringsforslag fremsat af Dansk Folkeparti (DF) er blevet forkastet, og lovforslaget er vedtaget uden fornyet udvalgsbehandling.'
$ws.Range("B120").Value = 'Behandlingen af statsborgerskab i Danmark rejser bekymringer om diskrimination, retssikkerhed og overholdelse af internationale konventioner, herunder handicapkonventionen, samtidig med at der er behov for en mere retfærdig og transparent proces for tildeling af statsborgerskab.'
$ws.Range("B121").Value = 'Dansk Folkeparti kritiserer regeringens statsborgerskabspolitik og mener, at tildeling bør baseres på respekt for danske værdier og samfundspositive bidrag.'
$ws.Range("B122").Value = 'Der rejses kritik af Folketingets udlændingepolitik, især vedrørende statsborgerskaber og integration, med fokus på behovet for ytringsfrihed, bekæmpelse af ekstremisme og en grundlæggende debat om danskhed og værdier.'
$ws.Range("B123").Value = 'Der er behov for en grundig politisk debat og stramninger i indfødsretslovgivningen for at sikre, at ansøgere opfylder danske værdier og undgå tildeling af statsborgerskab til personer med kriminel baggrund.'
$ws.Range("B124").Value = 'Debatten om dansk statsborgerskab centrerer sig om kravene til ansøgere, bekymringer om tildeling til personer med antidemokratiske holdninger, og behovet for en principiel tilgang, der sikrer, at nye borgere deler danske værdier.'
$ws.Range("B125").Value = 'Flere ændringsforslag fra DF blev forkastet, mens nogle forslag fra et flertal blev vedtaget, før lovforslaget endeligt blev godkendt.'
$ws.Range("B126").Value = 'Der er omfattende kritik af den nuværende administration af indfødsretsområdet, især vedrørende tildeling af dansk statsborgerskab til personer med kriminel baggrund, og der er enighed om behovet for strammere regler og bedre kontrol for at genoprette tilliden til systemet.'
$ws.Range("B127").Value = 'Diskussionen om indfødsret og statsborgerskab i Danmark er præget af afvisning af ændringsforslag fra Dansk Folkeparti, bekymringer om overholdelse af internationale konventioner, og en stram udlændingepolitik, der fokuserer på individuel vurdering af ansøgere.'
$ws.Range("B128").Value = 'Debatten om dansk statsborgerskab fokuserer på sprogkrav, kontrol af ansøgeres baggrund og holdninger, samt behovet for bedre screeningsmetoder og integration af personer uden statsborgerskab.'
$ws.Range("B129").Value = 'Debatten om lovforslaget om statsborgerskab i Danmark fokuserer på anerkendelse af de 1.100 ansøgere, stramme krav til tildeling, politiske uenigheder og bekymringer om integration og sikkerhed.'
$ws.Range("B130").Value = 'Der er en udbredt kritik af Folketingets udlændingepolitik med fokus på behovet for ændringer for at håndtere befolkningsudskiftning og bevare dansk kultur, samt en opfordring til politisk handling fra partier som Borgernes Parti og Dansk Folkeparti.'
